# Insert a new weekly price record for "Coliflor" (Agrícola del Norte S.A. de
# Arica) at row 22, pushing the existing rows (old rows 22-51) down by one
# (new rows 23-52). This mirrors how the data was edited upstream: a new
# daily/weekly observation was inserted into the middle of the historical
# series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22; everything below shifts down.
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = 44413
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100112008
$ws.Cells.Item(22, 7).Value = "Coliflor"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Tercera"
$ws.Cells.Item(22, 10).Value = 1500
$ws.Cells.Item(22, 11).Value = 500
$ws.Cells.Item(22, 12).Value = 600
$ws.Cells.Item(22, 13).Value = 550
$ws.Cells.Item(22, 14).Value = "`$/unidad"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 550
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
